{"js": "// 1. \"Office: M1.25\" -> \"Office: M1.30\"\nconst officeHits = context.document.body.search(\"Office: M1.25\", { matchCase: true });\nofficeHits.load(\"items\");\nawait context.sync();\n\nif (officeHits.items.length > 0) {\n  officeHits.items[0].insertText(\"Office: M1.30\", \"Replace\");\n  await context.sync();\n}\n\n// 2. \"Download and install Python on your own machine if you have one:\"\n//    -> \"Download and install Python (\" + bold(\"version 2.x\") + \") on your own machine if you have one:\"\nconst oldText = \"Download and install Python on your own machine if you have one:\";\nconst newText = \"Download and install Python (version 2.x) on your own machine if you have one:\";\n\nconst pyHits = context.document.body.search(oldText, { matchCase: true });\npyHits.load(\"items\");\nawait context.sync();\n\nif (pyHits.items.length > 0) {\n  // Replace the whole sentence first (plain text, single run).\n  pyHits.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n\n  // Then locate the \"version 2.x\" substring and make it bold \u2014 this splits\n  // the run into the correct before/bold/after pieces.\n  const boldHits = context.document.body.search(\"version 2.x\", { matchCase: true });\n  boldHits.load(\"items\");\n  await context.sync();\n\n  if (boldHits.items.length > 0) {\n    boldHits.items[0].font.bold = true;\n    await context.sync();\n  }\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. \"Office: M1.25\" -> \"Office: M1.30\"\n$officeRange = $d.Content\n$officeRange.Find.Execute(\"Office: M1.25\", $false, $false, $false, $false, $false, $true, 1, $false, \"Office: M1.30\", 2)\n\n# 2. \"Download and install Python on your own machine if you have one:\"\n#    -> \"Download and install Python (\" + bold(\"version 2.x\") + \") on your own machine if you have one:\"\n$oldText = \"Download and install Python on your own machine if you have one:\"\n$newText = \"Download and install Python (version 2.x) on your own machine if you have one:\"\n\n$pyRange = $d.Content\n$pyRange.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n# Find the newly inserted \"version 2.x\" substring and make it bold.\n$boldRange = $d.Content\n$found = $boldRange.Find.Execute(\"version 2.x\", $false, $false, $false, $false, $false, $true, 1, $false, \"\", 0)\nif ($found) {\n    $boldRange.Font.Bold = $true\n}\n"}
